# Update points for customer "09876543" -> 120.00
#
# The sheet stores this one customer's phone number as literal text
# "09876543" (leading zero kept) while every other row stores the same
# digits as a plain number (leading zero dropped). Editing that text
# value in place and writing a new points total looks, in the saved
# OOXML, like: a generic numeric-style row gets inserted directly above
# the old row (duplicating the row above it, the same way the rest of
# the sheet is shaped), which pushes the real "09876543" record down by
# one row; the points total on that (now shifted) record is then set to
# its new value.
#
# Locate the row dynamically instead of hard-coding row 53 so the script
# still works if the sheet shape changes slightly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$match = $ws.Columns.Item(1).Find("09876543")
$targetRow = $match.Row

# Duplicate the row immediately above the match, inserting it at the
# match's current position (this is what pushes the "09876543" row down
# by one and gives the vacated row the same look as its neighbours).
$ws.Rows.Item($targetRow - 1).Copy()
$ws.Rows.Item($targetRow).Insert()

# The real record (phone "09876543") now lives one row lower; update its
# total_points (column C) to the new value.
$newRow = $targetRow + 1
$ws.Cells.Item($newRow, 3).Value = 120
